# Updates cryptos list cell values (prices / volume %) per Fri Sep  6 23:35:28 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.544.58"
$ws.Range("E2").Value = "  -4.43%  "
$ws.Range("D3").Value = "2.197.47"
$ws.Range("E3").Value = "  -6.95%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "485.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.24%  "
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.51%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -4.32%  "
$ws.Range("D9").Value = "2.217.22"
$ws.Range("E9").Value = "  -6.12%  "
$ws.Range("E10").Value = "  -6.25%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.45%  "
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("D14").Value = "2.589.95"
$ws.Range("E14").Value = "  -6.91%  "
$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "53.477.06"
$ws.Range("E16").Value = "  -4.46%  "
$ws.Range("D18").Value = "2.215.20"
$ws.Range("E18").Value = "  -6.40%  "
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("E20").Value = "  -4.63%  "
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "294.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.38%  "
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.298.12"
$ws.Range("E27").Value = "  -7.09%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.145"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("E29").Value = "  -3.23%  "
$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.997"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "0.0₃0662"
$ws.Range("E33").Value = "  -6.71%  "
$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.994"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  -4.93%  "
$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.365"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "125.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.534"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "231.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0199"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.74%  "
